# Logged 2021 divisional round, simulated season from conference round
$wb = $excel.ActiveWorkbook
$wsRushing = $wb.Worksheets.Item("Rushing")
$wsReceiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet ---------------------------------------------------------
# Row 3 previously belonged to J.McNichols; D.Henry (back from injury for the
# divisional round) takes over as the Titans' featured rusher.
$wsRushing.Range("B3").Value = "D.Henry"
$wsRushing.Range("C3").Value = 127
$wsRushing.Range("D3").Value = 86
$wsRushing.Range("E3").Value = 34
$wsRushing.Range("F3").Value = 39

# Divisional-round stat adjustments for the rest of the rushing corps.
$wsRushing.Range("E2").Value = 16   # R.Tannehill 3DATT

$wsRushing.Range("C6").Value = 15   # D.Foreman 1DATT
$wsRushing.Range("D6").Value = 11   # D.Foreman 2DATT
$wsRushing.Range("E6").Value = 7    # D.Foreman 3DATT
$wsRushing.Range("F6").Value = 5    # D.Foreman RZATT

$wsRushing.Range("C7").Value = 1    # D.Hilliard 1DATT
$wsRushing.Range("D7").Value = 1    # D.Hilliard 2DATT
$wsRushing.Range("E7").Value = 1    # D.Hilliard 3DATT
$wsRushing.Range("F7").Value = 1    # D.Hilliard RZATT

# --- Receiving sheet --------------------------------------------------------
# A new row is added for D.Henry's receiving line; every other player shifts
# down one row and each player's season totals are refreshed for the simulated
# season through the conference round plus the divisional round game.

# Extend the table by one row, copying the formatting of the last data row
# (18) into the new row 19 so the appended row keeps consistent styling.
$wsReceiving.Range("A18:H18").Copy() | Out-Null
$wsReceiving.Range("A19:H19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$wsReceiving.Cells.Item(2, 1).Value = 0
$wsReceiving.Cells.Item(2, 2).Value = "D.Henry"
$wsReceiving.Cells.Item(2, 3).Value = 35
$wsReceiving.Cells.Item(2, 4).Value = 25
$wsReceiving.Cells.Item(2, 5).Value = 0
$wsReceiving.Cells.Item(2, 6).Value = 0
$wsReceiving.Cells.Item(2, 7).Value = 2
$wsReceiving.Cells.Item(2, 8).Value = 2

$wsReceiving.Cells.Item(3, 1).Value = 1
$wsReceiving.Cells.Item(3, 2).Value = "J.McNichols"
$wsReceiving.Cells.Item(3, 3).Value = 15
$wsReceiving.Cells.Item(3, 4).Value = 9
$wsReceiving.Cells.Item(3, 5).Value = 1
$wsReceiving.Cells.Item(3, 6).Value = 1
$wsReceiving.Cells.Item(3, 7).Value = 1
$wsReceiving.Cells.Item(3, 8).Value = 1

$wsReceiving.Cells.Item(4, 1).Value = 2
$wsReceiving.Cells.Item(4, 2).Value = "K.Blasingame"
$wsReceiving.Cells.Item(4, 3).Value = 3
$wsReceiving.Cells.Item(4, 4).Value = 2
$wsReceiving.Cells.Item(4, 5).Value = 0
$wsReceiving.Cells.Item(4, 6).Value = 0
$wsReceiving.Cells.Item(4, 7).Value = 0
$wsReceiving.Cells.Item(4, 8).Value = 0

$wsReceiving.Cells.Item(5, 1).Value = 3
$wsReceiving.Cells.Item(5, 2).Value = "D.Evans"
$wsReceiving.Cells.Item(5, 3).Value = 2
$wsReceiving.Cells.Item(5, 4).Value = 2
$wsReceiving.Cells.Item(5, 5).Value = 0
$wsReceiving.Cells.Item(5, 6).Value = 0
$wsReceiving.Cells.Item(5, 7).Value = 0
$wsReceiving.Cells.Item(5, 8).Value = 0

$wsReceiving.Cells.Item(6, 1).Value = 4
$wsReceiving.Cells.Item(6, 2).Value = "D.Foreman"
$wsReceiving.Cells.Item(6, 3).Value = 10
$wsReceiving.Cells.Item(6, 4).Value = 8
$wsReceiving.Cells.Item(6, 5).Value = 1
$wsReceiving.Cells.Item(6, 6).Value = 1
$wsReceiving.Cells.Item(6, 7).Value = 2
$wsReceiving.Cells.Item(6, 8).Value = 2

$wsReceiving.Cells.Item(7, 1).Value = 5
$wsReceiving.Cells.Item(7, 2).Value = "D.Hilliard"
$wsReceiving.Cells.Item(7, 3).Value = 3
$wsReceiving.Cells.Item(7, 4).Value = 3
$wsReceiving.Cells.Item(7, 5).Value = 0
$wsReceiving.Cells.Item(7, 6).Value = 0
$wsReceiving.Cells.Item(7, 7).Value = 0
$wsReceiving.Cells.Item(7, 8).Value = 0

$wsReceiving.Cells.Item(8, 1).Value = 6
$wsReceiving.Cells.Item(8, 2).Value = "A.Brown"
$wsReceiving.Cells.Item(8, 3).Value = 105
$wsReceiving.Cells.Item(8, 4).Value = 75
$wsReceiving.Cells.Item(8, 5).Value = 30
$wsReceiving.Cells.Item(8, 6).Value = 22
$wsReceiving.Cells.Item(8, 7).Value = 13
$wsReceiving.Cells.Item(8, 8).Value = 11

$wsReceiving.Cells.Item(9, 1).Value = 7
$wsReceiving.Cells.Item(9, 2).Value = "J.Jones"
$wsReceiving.Cells.Item(9, 3).Value = 61
$wsReceiving.Cells.Item(9, 4).Value = 46
$wsReceiving.Cells.Item(9, 5).Value = 11
$wsReceiving.Cells.Item(9, 6).Value = 8
$wsReceiving.Cells.Item(9, 7).Value = 7
$wsReceiving.Cells.Item(9, 8).Value = 5

$wsReceiving.Cells.Item(10, 1).Value = 8
$wsReceiving.Cells.Item(10, 2).Value = "C.Rogers"
$wsReceiving.Cells.Item(10, 3).Value = 33
$wsReceiving.Cells.Item(10, 4).Value = 23
$wsReceiving.Cells.Item(10, 5).Value = 4
$wsReceiving.Cells.Item(10, 6).Value = 2
$wsReceiving.Cells.Item(10, 7).Value = 4
$wsReceiving.Cells.Item(10, 8).Value = 1

$wsReceiving.Cells.Item(11, 1).Value = 9
$wsReceiving.Cells.Item(11, 2).Value = "R.McMath"
$wsReceiving.Cells.Item(11, 3).Value = 3
$wsReceiving.Cells.Item(11, 4).Value = 2
$wsReceiving.Cells.Item(11, 5).Value = 1
$wsReceiving.Cells.Item(11, 6).Value = 0
$wsReceiving.Cells.Item(11, 7).Value = 1
$wsReceiving.Cells.Item(11, 8).Value = 1

$wsReceiving.Cells.Item(12, 1).Value = 10
$wsReceiving.Cells.Item(12, 2).Value = "M.Johnson"
$wsReceiving.Cells.Item(12, 3).Value = 13
$wsReceiving.Cells.Item(12, 4).Value = 6
$wsReceiving.Cells.Item(12, 5).Value = 6
$wsReceiving.Cells.Item(12, 6).Value = 3
$wsReceiving.Cells.Item(12, 7).Value = 3
$wsReceiving.Cells.Item(12, 8).Value = 1

$wsReceiving.Cells.Item(13, 1).Value = 11
$wsReceiving.Cells.Item(13, 2).Value = "D.Fitzpatrick"
$wsReceiving.Cells.Item(13, 3).Value = 6
$wsReceiving.Cells.Item(13, 4).Value = 5
$wsReceiving.Cells.Item(13, 5).Value = 2
$wsReceiving.Cells.Item(13, 6).Value = 0
$wsReceiving.Cells.Item(13, 7).Value = 2
$wsReceiving.Cells.Item(13, 8).Value = 2

$wsReceiving.Cells.Item(14, 1).Value = 12
$wsReceiving.Cells.Item(14, 2).Value = "N.Westbrook-Ikhine"
$wsReceiving.Cells.Item(14, 3).Value = 28
$wsReceiving.Cells.Item(14, 4).Value = 21
$wsReceiving.Cells.Item(14, 5).Value = 11
$wsReceiving.Cells.Item(14, 6).Value = 5
$wsReceiving.Cells.Item(14, 7).Value = 5
$wsReceiving.Cells.Item(14, 8).Value = 5

$wsReceiving.Cells.Item(15, 1).Value = 13
$wsReceiving.Cells.Item(15, 2).Value = "C.Hollister"
$wsReceiving.Cells.Item(15, 3).Value = 6
$wsReceiving.Cells.Item(15, 4).Value = 4
$wsReceiving.Cells.Item(15, 5).Value = 1
$wsReceiving.Cells.Item(15, 6).Value = 0
$wsReceiving.Cells.Item(15, 7).Value = 2
$wsReceiving.Cells.Item(15, 8).Value = 1

$wsReceiving.Cells.Item(16, 1).Value = 14
$wsReceiving.Cells.Item(16, 2).Value = "A.Firkser"
$wsReceiving.Cells.Item(16, 3).Value = 36
$wsReceiving.Cells.Item(16, 4).Value = 30
$wsReceiving.Cells.Item(16, 5).Value = 3
$wsReceiving.Cells.Item(16, 6).Value = 1
$wsReceiving.Cells.Item(16, 7).Value = 5
$wsReceiving.Cells.Item(16, 8).Value = 4

$wsReceiving.Cells.Item(17, 1).Value = 15
$wsReceiving.Cells.Item(17, 2).Value = "M.Pruitt"
$wsReceiving.Cells.Item(17, 3).Value = 14
$wsReceiving.Cells.Item(17, 4).Value = 10
$wsReceiving.Cells.Item(17, 5).Value = 3
$wsReceiving.Cells.Item(17, 6).Value = 2
$wsReceiving.Cells.Item(17, 7).Value = 6
$wsReceiving.Cells.Item(17, 8).Value = 3

$wsReceiving.Cells.Item(18, 1).Value = 16
$wsReceiving.Cells.Item(18, 2).Value = "G.Swaim"
$wsReceiving.Cells.Item(18, 3).Value = 33
$wsReceiving.Cells.Item(18, 4).Value = 27
$wsReceiving.Cells.Item(18, 5).Value = 1
$wsReceiving.Cells.Item(18, 6).Value = 1
$wsReceiving.Cells.Item(18, 7).Value = 7
$wsReceiving.Cells.Item(18, 8).Value = 4

$wsReceiving.Cells.Item(19, 1).Value = 17
$wsReceiving.Cells.Item(19, 2).Value = "T.Hudson"
$wsReceiving.Cells.Item(19, 3).Value = 2
$wsReceiving.Cells.Item(19, 4).Value = 1
$wsReceiving.Cells.Item(19, 5).Value = 1
$wsReceiving.Cells.Item(19, 6).Value = 0
$wsReceiving.Cells.Item(19, 7).Value = 0
$wsReceiving.Cells.Item(19, 8).Value = 0

